$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.955.22'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.374.80'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('E7').Value = '  -2.39%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.624'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0932'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('E13').Value = '  -3.76%  '
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('E15').Value = '  -5.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.733.95'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.387.37'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.931.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '76.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.73'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '258.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.25%  '
$ws.Range('E24').Value = '  -3.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  -3.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('E29').Value = '  +3.03%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '172.14'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0901'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.74%  '
$ws.Range('E33').Value = '  +1.33%  '
$ws.Range('E34').Value = '  -6.79%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.121'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +12.64%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.132'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.68%  '
$ws.Range('E37').Value = '  -4.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0367'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.93'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.33%  '
$ws.Range('E40').Value = '  -4.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.54'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.243'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.04'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.68%  '
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.46'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.39'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '113.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.79%  '
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.28'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '77.45'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.76%  '
$ws.Range('E51').Value = '  -1.13%  '
